$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 1087.5
$ws.Range("I98").Value = 1497.5
$ws.Range("K98").Value = 1497.5
$ws.Range("M98").Value = 0.5
# Row 107
$ws.Range("H107").Value = 4913.6665
$ws.Range("J107").Value = 2493
$ws.Range("L107").Value = 2493
$ws.Range("N107").Value = -6333
# Row 122
$ws.Range("H122").Value = 1087.5
$ws.Range("I122").Value = 1497.5
$ws.Range("K122").Value = 4492.5
$ws.Range("M122").Value = -2042.5
# Row 137
$ws.Range("H137").Value = 1604
$ws.Range("I137").Value = 1277.8182
$ws.Range("K137").Value = 3833.4546
$ws.Range("M137").Value = -1283.4546
# Row 141
$ws.Range("H141").Value = 1017.94446
$ws.Range("I141").Value = 1026.4375
$ws.Range("K141").Value = 3079.3125
$ws.Range("M141").Value = 2100.6875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 911.2
$ws.Range("I61").Value = 913.55554
$ws.Range("K61").Value = 913.55554
$ws.Range("M61").Value = -701.55554
# Row 136
$ws.Range("H136").Value = 911.2
$ws.Range("I136").Value = 913.55554
$ws.Range("K136").Value = 2740.66662
$ws.Range("M136").Value = -190.66662

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 7441.3
$ws.Range("J3").Value = 10003.2
$ws.Range("L3").Value = 10003.2
$ws.Range("N3").Value = -10231.2
# Row 22
$ws.Range("H22").Value = 685
$ws.Range("I22").Value = 685
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 685
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -512
$ws.Range("N22").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1780.9
$ws.Range("I31").Value = 1318.3334
$ws.Range("K31").Value = 1318.3334
$ws.Range("M31").Value = -1023.3334
# Row 34
$ws.Range("H34").Value = 1780.9
$ws.Range("I34").Value = 1318.3334
$ws.Range("K34").Value = 1318.3334
$ws.Range("M34").Value = -1116.3334
# Row 58
$ws.Range("H58").Value = 1182.875
$ws.Range("I58").Value = 995.7143
$ws.Range("K58").Value = 995.7143
$ws.Range("M58").Value = -792.7143
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 94
$ws.Range("H94").Value = 4246.625
$ws.Range("I94").Value = 5040
$ws.Range("J94").Value = 3453.25
$ws.Range("K94").Value = 5040
$ws.Range("L94").Value = 3453.25
$ws.Range("M94").Value = -4589
$ws.Range("N94").Value = -4355.25
# Row 132
$ws.Range("H132").Value = 2499.44
$ws.Range("I132").Value = 2732.6365
$ws.Range("J132").Value = 789.3333
$ws.Range("K132").Value = 8197.9095
$ws.Range("L132").Value = 2367.9999
$ws.Range("M132").Value = -5667.9095
$ws.Range("N132").Value = -7427.9999
# Row 136
$ws.Range("H136").Value = 1182.875
$ws.Range("I136").Value = 995.7143
$ws.Range("K136").Value = 2987.1429
$ws.Range("M136").Value = -437.1428999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 28759.445
$ws.Range("I22").Value = 32260.625
$ws.Range("K22").Value = 96781.875
$ws.Range("M22").Value = -96612.875
# Row 27
$ws.Range("H27").Value = 28759.445
$ws.Range("I27").Value = 32260.625
$ws.Range("K27").Value = 96781.875
$ws.Range("M27").Value = -96679.875
# Row 80
$ws.Range("H80").Value = 8832.666999999999
$ws.Range("J80").Value = 10642.286
$ws.Range("L80").Value = 31926.858
$ws.Range("N80").Value = -33798.858
# Row 83
$ws.Range("H83").Value = 8832.666999999999
$ws.Range("J83").Value = 10642.286
$ws.Range("L83").Value = 95780.57399999999
$ws.Range("N83").Value = -105140.574

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4548.143
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 4548.143
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 80
$ws.Range("H80").Value = 2662.375
$ws.Range("J80").Value = 2862.5
$ws.Range("L80").Value = 2862.5
$ws.Range("N80").Value = -4858.5
# Row 83
$ws.Range("H83").Value = 2662.375
$ws.Range("J83").Value = 2862.5
$ws.Range("L83").Value = 14312.5
$ws.Range("N83").Value = -24296.5
# Row 113
$ws.Range("H113").Value = 2035.6
$ws.Range("I113").Value = 2046.5
$ws.Range("K113").Value = 2046.5
$ws.Range("M113").Value = 123.5
# Row 126
$ws.Range("H126").Value = 5597.5
$ws.Range("I126").Value = 4945
$ws.Range("K126").Value = 14835
$ws.Range("M126").Value = -12365

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3833
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3833
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3833
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -4423
# Row 27
$ws.Range("H27").Value = 3833
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3833
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3833
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -4047
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()
# Row 55
$ws.Range("H55").Value = 436.58334
$ws.Range("I55").Value = 116.333336
$ws.Range("J55").Value = 543.3333
$ws.Range("K55").Value = 116.333336
$ws.Range("L55").Value = 543.3333
$ws.Range("M55").Value = 56.666664
$ws.Range("N55").Value = -889.3333
# Row 132
$ws.Range("H132").Value = 2831.423
$ws.Range("J132").Value = 2897.4
$ws.Range("L132").Value = 8692.200000000001
$ws.Range("N132").Value = -13752.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4333
$ws.Range("I62").Value = 3849.5
$ws.Range("J62").Value = 5300
$ws.Range("K62").Value = 3849.5
$ws.Range("L62").Value = 5300
$ws.Range("M62").Value = -3225.5
$ws.Range("N62").Value = -6548
# Row 65
$ws.Range("H65").Value = 4333
$ws.Range("I65").Value = 3849.5
$ws.Range("J65").Value = 5300
$ws.Range("K65").Value = 19247.5
$ws.Range("L65").Value = 26500
$ws.Range("M65").Value = -16127.5
$ws.Range("N65").Value = -32740
